$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("donnees")

# Insert a new row before current row 7 ("% augmentation stats par level") to host "vit hero"
$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = "vit hero"
$ws.Range("B7").Value = 10

# Append new rows after the existing "% augmentation stats par level" row (now row 8)
$ws.Range("A9").Value = "pp attaque  1"
$ws.Range("B9").Value = 5

$ws.Range("A10").Value = "pp attaque  2"
$ws.Range("B10").Value = 6

$ws.Range("A11").Value = "pp attaque  3"
$ws.Range("B11").Value = 7

$ws.Range("A12").Value = "pp attaque  4"
$ws.Range("B12").Value = 8

$ws.Activate()
$ws.Range("E14").Select()
